$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 293.25
$ws.Range("I9").Value = 241
$ws.Range("K9").Value = 241
$ws.Range("M9").Value = -72

$ws.Range("H17").Value = 5883486
$ws.Range("J17").Value = 5883486
$ws.Range("L17").Value = 17650458
$ws.Range("N17").Value = -17650794

$ws.Range("H138").Value = 26318916
$ws.Range("I138").Value = 2669.6428
$ws.Range("J138").Value = 41670060
$ws.Range("K138").Value = 8008.928400000001
$ws.Range("L138").Value = 125010180
$ws.Range("M138").Value = -2868.928400000001
$ws.Range("N138").Value = -125020460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2361
$ws.Range("I2").Value = 2361
$ws.Range("K2").Value = 2361
$ws.Range("M2").Value = -2248

$ws.Range("H32").Value = 2623.98
$ws.Range("J32").Value = 664.5
$ws.Range("L32").Value = 664.5
$ws.Range("N32").Value = -1238.5

$ws.Range("H45").Value = 6444.4116
$ws.Range("J45").Value = 9999.333
$ws.Range("L45").Value = 9999.333
$ws.Range("N45").Value = -10753.333

$ws.Range("H61").Value = 17249368
$ws.Range("I61").Value = 25006696
$ws.Range("K61").Value = 25006696
$ws.Range("M61").Value = -25006484

$ws.Range("H116").Value = 2361
$ws.Range("I116").Value = 2361
$ws.Range("K116").Value = 2361
$ws.Range("M116").Value = -67

$ws.Range("H125").Value = 79996.75
$ws.Range("J125").Value = 79996.75
$ws.Range("L125").Value = 79996.75
$ws.Range("N125").Value = -89836.75

$ws.Range("H132").Value = 3290.9556
$ws.Range("I132").Value = 2746.3489
$ws.Range("K132").Value = 8239.046699999999
$ws.Range("M132").Value = -5709.046699999999

$ws.Range("H136").Value = 17249368
$ws.Range("I136").Value = 25006696
$ws.Range("K136").Value = 75020088
$ws.Range("M136").Value = -75017538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 5212
$ws.Range("I36").Value = 5212
$ws.Range("K36").Value = 5212
$ws.Range("M36").Value = -4678

$ws.Range("H86").Value = 2483.375
$ws.Range("I86").Value = 2462
$ws.Range("K86").Value = 2462
$ws.Range("M86").Value = -1339

$ws.Range("H89").Value = 2483.375
$ws.Range("I89").Value = 2462
$ws.Range("K89").Value = 12310
$ws.Range("M89").Value = -6694

$ws.Range("H105").Value = 3304.7
$ws.Range("I105").Value = 2906.8333
$ws.Range("K105").Value = 2906.8333
$ws.Range("M105").Value = -1159.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4126.606
$ws.Range("I31").Value = 3188.6191
$ws.Range("K31").Value = 3188.6191
$ws.Range("M31").Value = -2893.6191

$ws.Range("H34").Value = 4126.606
$ws.Range("I34").Value = 3188.6191
$ws.Range("K34").Value = 3188.6191
$ws.Range("M34").Value = -2986.6191

$ws.Range("H41").Value = 22156.111
$ws.Range("J41").Value = 22667.834
$ws.Range("L41").Value = 22667.834
$ws.Range("N41").Value = -23523.834

$ws.Range("H51").Value = 24097.5
$ws.Range("J51").Value = 28796.666
$ws.Range("L51").Value = 28796.666
$ws.Range("N51").Value = -30268.666

$ws.Range("H58").Value = 6188.0454
$ws.Range("I58").Value = 2418.625
$ws.Range("K58").Value = 2418.625
$ws.Range("M58").Value = -2215.625

$ws.Range("H61").Value = 24097.5
$ws.Range("J61").Value = 28796.666
$ws.Range("L61").Value = 28796.666
$ws.Range("N61").Value = -29492.666

$ws.Range("H86").Value = 7668.4287
$ws.Range("I86").Value = 6295
$ws.Range("J86").Value = 9499.667
$ws.Range("K86").Value = 6295
$ws.Range("L86").Value = 9499.667
$ws.Range("M86").Value = -5172
$ws.Range("N86").Value = -11745.667

$ws.Range("H89").Value = 7668.4287
$ws.Range("I89").Value = 6295
$ws.Range("J89").Value = 9499.667
$ws.Range("K89").Value = 31475
$ws.Range("L89").Value = 47498.335
$ws.Range("M89").Value = -25859
$ws.Range("N89").Value = -58730.335

$ws.Range("H105").Value = 1983.3636
$ws.Range("I105").Value = 2190.7778
$ws.Range("K105").Value = 2190.7778
$ws.Range("M105").Value = -443.7777999999998

$ws.Range("H136").Value = 6188.0454
$ws.Range("I136").Value = 2418.625
$ws.Range("K136").Value = 7255.875
$ws.Range("M136").Value = -4705.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 983.6667
$ws.Range("I8").Value = 983.6667
$ws.Range("K8").Value = 2951.0001
$ws.Range("M8").Value = -2812.0001

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H122").Value = 3174.625
$ws.Range("I122").Value = 20000
$ws.Range("J122").Value = 771
$ws.Range("K122").Value = 180000
$ws.Range("L122").Value = 6939
$ws.Range("M122").Value = -177550
$ws.Range("N122").Value = -11839

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H36").Value = 3224.75
$ws.Range("I36").Value = 2199.5
$ws.Range("J36").Value = 4250
$ws.Range("K36").Value = 2199.5
$ws.Range("L36").Value = 4250
$ws.Range("M36").Value = -1714.5
$ws.Range("N36").Value = -5220

$ws.Range("H122").Value = 3942.9285
$ws.Range("I122").Value = 2291.2727
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 6873.8181
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -4423.8181
$ws.Range("N122").Value = -34897

$ws.Range("H132").Value = 3536.15
$ws.Range("I132").Value = 3031.0322
$ws.Range("K132").Value = 9093.0966
$ws.Range("M132").Value = -6563.096600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1400
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -5

$ws.Range("H27").Value = 1400
$ws.Range("I27").Value = 300
$ws.Range("K27").Value = 300
$ws.Range("M27").Value = -193

$ws.Range("H46").Value = 21743.666
$ws.Range("J46").Value = 21743.666
$ws.Range("L46").Value = 21743.666
$ws.Range("N46").Value = -22119.666

$ws.Range("H63").Value = 49934.6
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 49934.6
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 49934.6
$ws.Range("N63").Value = -51432.6
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 49934.6
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 49934.6
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 149803.8
$ws.Range("N66").Value = -157291.8
$ws.Range("M66").ClearContents()

$ws.Range("H132").Value = 7489.5312
$ws.Range("I132").Value = 13572.462
$ws.Range("K132").Value = 40717.386
$ws.Range("M132").Value = -38187.386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 2022
$ws.Range("I43").Value = 27.75
$ws.Range("K43").Value = 27.75
$ws.Range("M43").Value = 121.25

$ws.Range("H68").Value = 39000
$ws.Range("J68").Value = 39000
$ws.Range("L68").Value = 39000
$ws.Range("N68").Value = -40622

$ws.Range("H71").Value = 39000
$ws.Range("J71").Value = 39000
$ws.Range("L71").Value = 117000
$ws.Range("N71").Value = -125112
